$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feature Injection")

# ---------------------------------------------------------------------------
# The "Posting Label" field table lives in columns B (field name) and C
# (field value), rows 2..8. We are inserting a brand-new field
# ("data.kind" = "capability-hierarchy") right after the "excelAPI" row
# (row 2), which pushes every subsequent field down by one row. The very
# last field (old "dataRange" = "E5:I100") ends up in the newly-created
# row 9, and is itself renamed to "data.range".
#
# Columns E..Q (the actual BDD Feature-Injection content/grid) must NOT be
# shifted - they keep referring to the same row numbers as before.
# ---------------------------------------------------------------------------

# 1) Snapshot the current B/C values (rows 3..8) before we overwrite anything.
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$b4 = $ws.Range("B4").Value2
$c4 = $ws.Range("C4").Value2
$b5 = $ws.Range("B5").Value2
$c5 = $ws.Range("C5").Value2
$b6 = $ws.Range("B6").Value2
$c6 = $ws.Range("C6").Value2
$b7 = $ws.Range("B7").Value2
$c7 = $ws.Range("C7").Value2
$b8 = $ws.Range("B8").Value2
$c8 = $ws.Range("C8").Value2

# 2) Give row 9 the same cell formatting as the other filled-in field rows
#    (it used to be a blank placeholder row) before we put real values in it.
$ws.Range("B8:C8").Copy() | Out-Null
$ws.Range("B9:C9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Shift the field rows down by one (bottom to top so nothing is clobbered
#    before it has been read).
$ws.Range("B9").Value2 = $b8
$ws.Range("C9").Value2 = $c8
$ws.Range("B8").Value2 = $b7
$ws.Range("C8").Value2 = $c7
$ws.Range("B7").Value2 = $b6
$ws.Range("C7").Value2 = $c6
$ws.Range("B6").Value2 = $b5
$ws.Range("C6").Value2 = $c5
$ws.Range("B5").Value2 = $b4
$ws.Range("C5").Value2 = $c4
$ws.Range("B4").Value2 = $b3
$ws.Range("C4").Value2 = $c3

# 4) Update the posting label value that changed (row 2's "excelAPI" now
#    points at the BDD manifest, not the capability-hierarchy one).
$ws.Range("C2").Value2 = "bdd.kernel.a6i.xlsx/v1a"

# 5) Put the new field into row 3.
$ws.Range("C3").Value2 = "capability-hierarchy"
$ws.Range("B3").Value2 = "data.kind"

# 6) Rename the old "dataRange" field (now living in row 9) to "data.range".
$ws.Range("B9").Value2 = "data.range"

# ---------------------------------------------------------------------------
# Row heights: the field row that wraps text onto two lines ("BDD tests")
# moved from row 3 to row 4, and the explicit row height followed it.
# ---------------------------------------------------------------------------
$ws.Rows(4).RowHeight = 28.5
$ws.Rows(3).EntireRow.AutoFit()

# ---------------------------------------------------------------------------
# Selection: the saved cell selection moved to C16.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C16").Select()
